# Shopping / merchant functionality: new Dialog rows, a "price" column on
# Items, a "gold" column on Individuals, and new trigger/event columns +
# a merchant "attack" event row on Events.

$wb = $excel.ActiveWorkbook

$wsDialog      = $wb.Worksheets.Item("Dialog")
$wsItems       = $wb.Worksheets.Item("Items")
$wsIndividuals = $wb.Worksheets.Item("Individuals")
$wsEvents      = $wb.Worksheets.Item("Events")

# ---------------------------------------------------------------------
# Write the brand-new pieces of text first, and in the same order the
# author originally typed them in, so new shared-string entries line up
# the same way they did in the real edit (Events trigger headers, then
# the merchant attack line, then the two new dialog lines, then the new
# "price" and "gold" column headers).
# ---------------------------------------------------------------------
$wsEvents.Range("J1").Value = "Triggers:"
$wsEvents.Range("K1").Value = "Attack"
$wsEvents.Range("N1").Value = "Harm"
$wsEvents.Range("Q1").Value = "Death"
$wsEvents.Range("L2").Value = "EventID"
$wsEvents.Range("H3").Value = '"Have at you!"'

$wsDialog.Range("B11").Value = "Welcome! We've got the lowest prices!"
$wsDialog.Range("B12").Value = "What would you like?"

$wsItems.Range("H1:H7").EntireColumn.Insert()
$wsItems.Range("H1").Value = "price"

$wsIndividuals.Range("AG1").Value = "gold"

# ---------------------------------------------------------------------
# Sheet "Dialog" (sheet1): finish the new merchant dialog entries
# (IDs 1004, 1005) and four checkpoint rows referencing dialog 1005.
# ---------------------------------------------------------------------
$wsDialog.Range("A11").Value = 1004
$wsDialog.Range("C11").Value = 0
$wsDialog.Range("D11").Value = 1005
$wsDialog.Range("E11").Value = 0
$wsDialog.Range("F11").Value = 0

$wsDialog.Range("A12").Value = 1005
$wsDialog.Range("D12").Value = 0
$wsDialog.Range("E12").Value = 0
$wsDialog.Range("F12").Value = 0

$wsDialog.Range("A17").Value = 1005
$wsDialog.Range("A18").Value = 1005
$wsDialog.Range("A19").Value = 1005
$wsDialog.Range("A20").Value = 1005

# ---------------------------------------------------------------------
# Sheet "Items" (sheet2): the "price" column was inserted above right
# after itemType (column G); fill in the per-item price values.
# ---------------------------------------------------------------------
$wsItems.Range("H2").Value = 25
$wsItems.Range("H3").Value = 10
$wsItems.Range("H4").Value = 20
$wsItems.Range("H5").Value = 15
$wsItems.Range("H6").Value = 25
$wsItems.Range("H7").Value = 50

# ---------------------------------------------------------------------
# Sheet "Individuals" (sheet3): fill in the "gold" column after the
# existing dialogID column (AF), carrying how much gold each
# individual/merchant has.
# ---------------------------------------------------------------------
$wsIndividuals.Range("AG2").Value = 15
$wsIndividuals.Range("AG3").Value = 5
$wsIndividuals.Range("AG4").Value = 5
$wsIndividuals.Range("AG5").Value = 5
$wsIndividuals.Range("AG6").Value = 5
$wsIndividuals.Range("AG7").Value = 5
$wsIndividuals.Range("AG8").Value = 30

# ---------------------------------------------------------------------
# Sheet "Events" (sheet4): finish the Attack/Harm/Death trigger columns
# (IndividualID + EventID pairs) and the sample Attack trigger row
# pointing at the merchant (IndividualID 206) and its "attack" event.
# ---------------------------------------------------------------------
$wsEvents.Range("K2").Value = "IndividualID"
$wsEvents.Range("N2").Value = "IndividualID"
$wsEvents.Range("O2").Value = "EventID"
$wsEvents.Range("Q2").Value = "IndividualID"
$wsEvents.Range("R2").Value = "EventID"

$wsEvents.Range("K3").Value = 206
$wsEvents.Range("L3").Value = 2

$wsEvents.Columns.Item(11).ColumnWidth = 11.7109375
$wsEvents.Columns.Item(14).ColumnWidth = 11.7109375
$wsEvents.Columns.Item(17).ColumnWidth = 11.7109375

# ---------------------------------------------------------------------
# View state: restore the selections the author left on each sheet and
# make "Individuals" the active tab (matches activeTab going 3 -> 2).
# ---------------------------------------------------------------------
$wsDialog.Range("B17").Select()
$wsItems.Range("A2:AO7").Select()
$wsEvents.Range("K5").Select()
$wsIndividuals.Range("AB10").Select()
$wsIndividuals.Activate()

Write-Host "Merchant/shopping edit applied."
